$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "price" text (e.g. "69.167.63", "602.29") that must stay
# plain text, same as the source inline strings -- not get reinterpreted as
# numbers by Excel. Force each target cell to Text format individually
# (a multi-area union range only honors NumberFormat on its first area).
$priceCells = @("D2","D3","D5","D6","D11","D13","D15","D16","D17","D19","D21","D22","D23","D25","D27","D28","D31","D33","D34","D35","D38","D41","D42","D43","D44","D45","D46","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.167.63"
$ws.Range("D3").Value = "3.749.45"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "602.29"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "166.98"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "38.04"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "4.379.22"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "3.750.06"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "69.172.47"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "17.40"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  +18.74%  "
$ws.Range("D22").Value = "493.71"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").Value = "84.81"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "2.48"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "31.58"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "3.895.62"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "3.687.85"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "5.99"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").Value = "0.325"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "3.00"
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D43").Value = "48.85"
$ws.Range("D44").Value = "429.28"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "1.99"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "8.47"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D48").Value = "40.20"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "141.07"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "2.796.83"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("E51").Value = "  +0.04%  "
